$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Add()
$ws2.Name = "Planilha1"

$r = $ws2.Range("B3")
$r.Value = '> dm_tests[["pvalues"]]'
$r.Font.Name = "Cascadia Code"
$r.Font.Size = 7
$r.Font.Color = 0xA85989
$r.VerticalAlignment = -4108

$r2 = $ws2.Range("B8")
$r2.Value = '> dm_test_mean = compute_dmv2()'
$r2.Font.Name = "Cascadia Code"
$r2.Font.Size = 7
$r2.Font.Color = 0xA85989
$r2.VerticalAlignment = -4108
$r2.Interior.Color = 0xFFFFFF

$ws2.Move($null, $wb.Worksheets.Item("Sheet1"))
